$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores every value as plain text in the source
# workbook (even purely-numeric-looking strings such as "1.00"), so
# that values like "2.40" or "304.89" keep their exact original
# formatting instead of becoming floating point numbers. Pre-format
# the updated Price cells whose new text is otherwise a valid number
# as Text ("@") before assigning, so COM's Value setter keeps the
# literal digits/trailing zeros instead of auto-converting.

$ws.Range("D2").Value = '42.731.53'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '2.295.32'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.89'
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.31'
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").Value = '  -1.65%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.07'
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.64'
$ws.Range("E12").Value = '  +5.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.119'
$ws.Range("E13").Value = '  +2.02%  '
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").Value = '2.653.23'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '2.291.60'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '42.656.27'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.09'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.57'
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.60'
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '166.26'
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.84'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.95'
$ws.Range("E33").Value = '  +5.05%  '
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.43'
$ws.Range("E35").Value = '  -7.05%  '
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.101'
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("E40").Value = '  -0.52%  '
$ws.Range("E41").Value = '  -1.60%  '
$ws.Range("D42").Value = '1.998.87'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.23'
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.11'
$ws.Range("E45").Value = '  +5.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.12'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.75'
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.43'
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.520.85'
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.82'
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.90'
$ws.Range("E51").Value = '  -1.17%  '
